# Lab1Rubric_CS295N.xlsx - "Updated the lab3 assignment" commit
#
# The old single-sheet "Sheet1" (rubric w/ Possible+Actual columns) is
# renamed to "Rubric" and turned into a blank template (Possible column
# only). A new "Grade" sheet is added after it, holding the filled-in
# rubric (Possible + Actual columns) that used to live on Sheet1 - this
# becomes the active/selected sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the original sheet and add the new one right after it.
# ---------------------------------------------------------------------
$rubric = $wb.Worksheets.Item(1)
$rubric.Name = "Rubric"

$grade = $wb.Worksheets.Add($null, $rubric)
$grade.Name = "Grade"

# ---------------------------------------------------------------------
# 2. Rebuild "Rubric" (template: Requirements / Possible only).
# ---------------------------------------------------------------------
$rubric.Cells.Clear()

$rubric.Range("A1").Value = "Lab 1"
$rubric.Range("A1").Font.Bold = $true

$rubric.Range("A3").Value = "Requirements"
$rubric.Range("A3").Font.Underline = $true
$rubric.Range("B3").Value = "Possible"
$rubric.Range("B3").Font.Underline = $true
$rubric.Range("C3").Font.Underline = $true

$rubric.Range("A4").Value = "MVC site"

$rubric.Range("A5").Value = "Project builds and runs"
$rubric.Range("A5").Font.Italic = $true
$rubric.Range("B5").Value = 20

$rubric.Range("A6").Value = "Site name"
$rubric.Range("A6").Font.Italic = $true
$rubric.Range("B6").Value = 5

$rubric.Range("A7").Value = "Home page title and text"
$rubric.Range("A7").Font.Italic = $true
$rubric.Range("B7").Value = 5

$rubric.Range("A8").Font.Italic = $true

$rubric.Range("A9").Value = "GitHub repository"

$rubric.Range("A10").Value = "main branch"
$rubric.Range("A10").Font.Italic = $true
$rubric.Range("B10").Value = 5

$rubric.Range("A11").Value = ".gitignore"
$rubric.Range("A11").Font.Italic = $true
$rubric.Range("B11").Value = 5

$rubric.Range("A12").Value = "lab1 branch"
$rubric.Range("A12").Font.Italic = $true
$rubric.Range("B12").Value = 5

$rubric.Range("A13").Value = "Instructor invited"
$rubric.Range("A13").Font.Italic = $true
$rubric.Range("B13").Value = 5

$rubric.Range("A15").Value = "Total"
$rubric.Range("A15").Font.Italic = $true
$rubric.Range("B15").Formula = "=SUM(B4:B13)"

$rubric.Columns.Item(1).ColumnWidth = 24.33

$rubric.Range("E16").Select() | Out-Null

# ---------------------------------------------------------------------
# 3. Build "Grade" (filled-in: Requirements / Possible / Actual).
# ---------------------------------------------------------------------
$grade.Range("A1").Value = "Lab 1"
$grade.Range("A1").Font.Bold = $true

$grade.Range("A2").Value = "Excellent work! Everything looks great."

$grade.Range("A3").Value = "Requirements"
$grade.Range("A3").Font.Underline = $true
$grade.Range("B3").Value = "Possible"
$grade.Range("B3").Font.Underline = $true
$grade.Range("C3").Value = "Actual"
$grade.Range("C3").Font.Underline = $true

$grade.Range("A4").Value = "MVC site"

$grade.Range("A5").Value = "Project builds and runs"
$grade.Range("A5").Font.Italic = $true
$grade.Range("B5").Value = 20
$grade.Range("C5").Value = 20

$grade.Range("A6").Value = "Site name"
$grade.Range("A6").Font.Italic = $true
$grade.Range("B6").Value = 5
$grade.Range("C6").Value = 5

$grade.Range("A7").Value = "Home page title and text"
$grade.Range("A7").Font.Italic = $true
$grade.Range("B7").Value = 5
$grade.Range("C7").Value = 5

$grade.Range("A8").Font.Italic = $true

$grade.Range("A9").Value = "GitHub repository"

$grade.Range("A10").Value = "main branch"
$grade.Range("A10").Font.Italic = $true
$grade.Range("B10").Value = 5
$grade.Range("C10").Value = 5

$grade.Range("A11").Value = ".gitignore"
$grade.Range("A11").Font.Italic = $true
$grade.Range("B11").Value = 5
$grade.Range("C11").Value = 5

$grade.Range("A12").Value = "lab1 branch"
$grade.Range("A12").Font.Italic = $true
$grade.Range("B12").Value = 5
$grade.Range("C12").Value = 5

$grade.Range("A13").Value = "Instructor invited"
$grade.Range("A13").Font.Italic = $true
$grade.Range("B13").Value = 5
$grade.Range("C13").Value = 5

$grade.Range("A15").Value = "Total"
$grade.Range("A15").Font.Bold = $true
$grade.Range("B15").Formula = "=SUM(B4:B13)"
$grade.Range("B15").Font.Bold = $true
$grade.Range("C15").Formula = "=SUM(C4:C13)"
$grade.Range("C15").Font.Bold = $true

$grade.Range("A2:C15").Select() | Out-Null

# ---------------------------------------------------------------------
# 4. "Grade" is the tab the author left active.
# ---------------------------------------------------------------------
$grade.Activate() | Out-Null
